$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New label above existing table ---
$ws.Cells.Item(4, 5).Value = "WITH YEAR"

# --- New "WITHOUT YEAR" table (rows 14-19) ---
$ws.Cells.Item(14, 5).Value = "WITHOUT YEAR"

# Header row (15)
$ws.Cells.Item(15, 5).Value = "Data Type"
$ws.Cells.Item(15, 6).Value = "Entire_Scale"
$ws.Cells.Item(15, 7).Value = "Growing_Season"
$ws.Cells.Item(15, 8).Value = "Every_4_months"
$ws.Cells.Item(15, 9).Value = "Every_2_months "
$ws.Cells.Item(15, 10).Value = "Monthly "

# Row 16 - R^2_Train
$ws.Cells.Item(16, 5).Value = "R^2_Train"
$ws.Cells.Item(16, 6).Value = 0.93095170000000005
$ws.Cells.Item(16, 7).Value = 0.89611419999999997
$ws.Cells.Item(16, 8).Value = 0.91637990000000002
$ws.Cells.Item(16, 9).Value = 0.92126090000000005
$ws.Cells.Item(16, 10).Value = 0.9421235

# Row 17 - R^2_Test
$ws.Cells.Item(17, 5).Value = "R^2_Test"
$ws.Cells.Item(17, 6).Value = 0.26426240000000001
$ws.Cells.Item(17, 7).Value = 0.1032163
$ws.Cells.Item(17, 8).Value = 0.13765759999999999
$ws.Cells.Item(17, 9).Value = 0.2008228
$ws.Cells.Item(17, 10).Value = 0.24345849999999999

# Row 18 - RSME_Train
$ws.Cells.Item(18, 5).Value = "RSME_Train"
$ws.Cells.Item(18, 6).Value = 1.1101209999999999
$ws.Cells.Item(18, 7).Value = 1.2833749999999999
$ws.Cells.Item(18, 8).Value = 1.1861950000000001
$ws.Cells.Item(18, 9).Value = 1.158674
$ws.Cells.Item(18, 10).Value = 1.0812889999999999

# Row 19 - RSME_Test
$ws.Cells.Item(19, 5).Value = "RSME_Test"
$ws.Cells.Item(19, 6).Value = 2.8240430000000001
$ws.Cells.Item(19, 7).Value = 3.1396310000000001
$ws.Cells.Item(19, 8).Value = 3.0429409999999999
$ws.Cells.Item(19, 9).Value = 2.9239259999999998
$ws.Cells.Item(19, 10).Value = 2.8593310000000001

# --- Update selection to match the author's final cursor position ---
$ws.Range("L15:M19").Select()
